$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_suite")

# Flip the "Runmode" flag from "Y" to "N" for the Ordering, Parameter,
# VerifyLoginPage and Registration test cases (Login stays "Y").
$ws.Range("B3:B6").Value = "N"

# Move the active selection to B18 (single cell) instead of B2:B6.
$ws.Range("B18").Select()
